# Append-scrape update: 2025-11-30 01:29 JST
#
# Sheet "ランサーズ" (sheet 1) is refreshed with 3 newly scraped job
# postings merged into the existing 8, re-sorted by "優先度スコア"
# (column G) descending, and every row's "取得日時" timestamp bumped to
# the new scrape time.
#
# The combined/sorted list (score desc):
#   360 AI-Gravity study-group talk      (existing, row 2, unchanged pos)
#   323 React/TypeScript frontend        (existing, row 3, unchanged pos)
#   303 BlockChain/Solidity programmer   (NEW)                -> row 4
#   170 Inventory/sales mgmt tool        (existing, was row 4) -> row 5
#    68 Pet EC platform engineer         (existing, was row 5) -> row 6
#    60 Mansion mgmt-union system        (existing, was row 6) -> row 7
#    45 Appealing WEB site freelancer    (NEW)                -> row 8
#    38 Apache Answer Q&A server setup   (existing, was row 7) -> row 9
#    33 WordPress render-blocking fix    (NEW)                -> row 10
#    13 Wartales weapon icon swap        (existing, was row 8) -> row 11
#    10 Local subsidy expert wanted      (existing, was row 9) -> row 12
#
# The first 8 rows (2-9) are updated in place (plain value overwrites), so
# the worksheet's existing <hyperlinks> entries for F2:F9 are left exactly
# as-is (same r:id's, now pointing at whichever URL used to live in that
# slot) -- matching the source workbook's own scraper behaviour. Only the
# 3 genuinely new rows appended at the bottom (10-12) receive fresh
# hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$ts = "2025-11-30 01:29:30"

# --- Rows 2-3: timestamp refresh only, all other columns unchanged --------
$ws.Cells.Item(2, 1).Value = $ts
$ws.Cells.Item(3, 1).Value = $ts

# --- Row 4: now the BlockChain/Solidity posting (NEW) ----------------------
$ws.Cells.Item(4, 1).Value = $ts
$ws.Cells.Item(4, 2).Value = "【急募】BlockChainとSolidityに精通したプログラマー募集"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5443998"
$ws.Cells.Item(4, 7).Value = 303
$ws.Cells.Item(4, 8).Value = "🔥AI,Ai"

# --- Row 5: now holds the old row-4 posting (在庫・販売管理ツール) ---------
$ws.Cells.Item(5, 1).Value = $ts
$ws.Cells.Item(5, 2).Value = "【急募】在庫・販売管理ツールの開発依頼"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5443889"
$ws.Cells.Item(5, 7).Value = 170
$ws.Cells.Item(5, 8).Value = "◆ツール,開発 ◇管理"

# --- Row 6: now holds the old row-5 posting (ペット向けEC) -----------------
$ws.Cells.Item(6, 1).Value = $ts
$ws.Cells.Item(6, 2).Value = "【急募】革新的ペット向けECプラットフォーム開発エンジニア募集"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5443928"
$ws.Cells.Item(6, 7).Value = 68
$ws.Cells.Item(6, 8).Value = "◆開発"

# --- Row 7: now holds the old row-6 posting (マンション管理組合) ----------
$ws.Cells.Item(7, 1).Value = $ts
$ws.Cells.Item(7, 2).Value = "マンション管理組合のシステム設計構築依頼"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5443592"
$ws.Cells.Item(7, 7).Value = 60
$ws.Cells.Item(7, 8).Value = "◇管理"

# --- Row 8: now the WEB site freelancer posting (NEW) -----------------------
$ws.Cells.Item(8, 1).Value = $ts
$ws.Cells.Item(8, 2).Value = "【急募】魅力的なWEBサイト制作のフリーランスを探しています!"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5444036"
$ws.Cells.Item(8, 7).Value = 45
$ws.Cells.Item(8, 8).Value = "◇サイト"

# --- Row 9: now holds the old row-7 posting (Apache Answer構築) -----------
$ws.Cells.Item(9, 1).Value = $ts
$ws.Cells.Item(9, 2).Value = "【Apache Answer構築】弁護士ドットコムのような専門家Q&Aサイトのサーバー構築・初期設定"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5443617"
$ws.Cells.Item(9, 7).Value = 38
$ws.Cells.Item(9, 8).Value = "◇サイト"

# --- Row 10: brand-new appended row (WordPress posting) --------------------
$ws.Cells.Item(10, 1).Value = $ts
$ws.Cells.Item(10, 2).Value = "wordpressレンダリングを妨げるリソースの除外"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5016989"
$ws.Cells.Item(10, 7).Value = 33
$ws.Cells.Item(10, 8).Value = "○WordPress"

# --- Row 11: brand-new appended row (old row-8 posting, Wartales) ----------
$ws.Cells.Item(11, 1).Value = $ts
$ws.Cells.Item(11, 2).Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5443568"
$ws.Cells.Item(11, 7).Value = 13

# --- Row 12: brand-new appended row (old row-9 posting, 補助金) ------------
$ws.Cells.Item(12, 1).Value = $ts
$ws.Cells.Item(12, 2).Value = "地方の補助金に詳しい方募集"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5443921"
$ws.Cells.Item(12, 7).Value = 10

# --- New hyperlinks only for the 3 newly appended rows ----------------------
# (F2:F9 already carry the correct, untouched hyperlink relationships.)
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), "https://www.lancers.jp/work/detail/5016989")
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), "https://www.lancers.jp/work/detail/5443568")
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), "https://www.lancers.jp/work/detail/5443921")

# Re-apply the plain "Hyperlink" cell style (same style index the F2:F9
# hyperlink cells already use) so the new cells match the sheet's existing
# hyperlink formatting instead of getting a freshly-minted duplicate style.
$ws.Cells.Item(10, 6).Style = "Hyperlink"
$ws.Cells.Item(11, 6).Style = "Hyperlink"
$ws.Cells.Item(12, 6).Style = "Hyperlink"
